$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add value to F19 (new data point)
$ws.Range("F19").Value = 20

# Update selection to F20
$ws.Range("F20").Select()
